$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 3300
$ws.Range("I43").Value = 5466.6665
$ws.Range("J43").Value = 1675
$ws.Range("K43").Value = 5466.6665
$ws.Range("L43").Value = 1675
$ws.Range("M43").Value = -5397.6665
$ws.Range("N43").Value = -1813
# Row 98
$ws.Range("H98").Value = 2961.4285
$ws.Range("I98").Value = 2875.3333
$ws.Range("K98").Value = 2875.3333
$ws.Range("M98").Value = -1377.3333
# Row 112
$ws.Range("H112").Value = 3242.8572
$ws.Range("J112").Value = 3340.625
$ws.Range("L112").Value = 10021.875
$ws.Range("N112").Value = -12237.875
# Row 113
$ws.Range("H113").Value = 2924.8438
$ws.Range("I113").Value = 2780.3845
$ws.Range("J113").Value = 3550.8333
$ws.Range("K113").Value = 2780.3845
$ws.Range("L113").Value = 3550.8333
$ws.Range("M113").Value = 473.6154999999999
$ws.Range("N113").Value = -10058.8333
# Row 116
$ws.Range("H116").Value = 4523.077
$ws.Range("I116").Value = 5033.3335
$ws.Range("J116").Value = 3375
$ws.Range("K116").Value = 5033.3335
$ws.Range("L116").Value = 3375
$ws.Range("M116").Value = -1591.3335
$ws.Range("N116").Value = -10259
# Row 122
$ws.Range("H122").Value = 2961.4285
$ws.Range("I122").Value = 2875.3333
$ws.Range("K122").Value = 8625.999899999999
$ws.Range("M122").Value = -6175.999899999999
# Row 132
$ws.Range("H132").Value = 2327773.2
$ws.Range("I132").Value = 2207.7812
$ws.Range("J132").Value = 9093055
$ws.Range("K132").Value = 6623.3436
$ws.Range("L132").Value = 27279165
$ws.Range("M132").Value = -4093.3436
$ws.Range("N132").Value = -27284225
# Row 137
$ws.Range("H137").Value = 1630.3334
$ws.Range("I137").Value = 1245.5
$ws.Range("J137").Value = 2400
$ws.Range("K137").Value = 3736.5
$ws.Range("L137").Value = 7200
$ws.Range("M137").Value = -1186.5
$ws.Range("N137").Value = -12300
# Row 138
$ws.Range("H138").Value = 4034.237
$ws.Range("I138").Value = 2690.1177
$ws.Range("J138").Value = 4319.8623
$ws.Range("K138").Value = 8070.353099999999
$ws.Range("L138").Value = 12959.5869
$ws.Range("M138").Value = -2930.353099999999
$ws.Range("N138").Value = -23239.5869

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 18951.879
$ws.Range("I32").Value = 13614.026
$ws.Range("J32").Value = 48310.07
$ws.Range("K32").Value = 13614.026
$ws.Range("L32").Value = 48310.07
$ws.Range("M32").Value = -13327.026
$ws.Range("N32").Value = -48884.07
# Row 74
$ws.Range("H74").Value = 941.36957
$ws.Range("I74").Value = 889.87177
$ws.Range("J74").Value = 1228.2858
$ws.Range("K74").Value = 889.87177
$ws.Range("L74").Value = 1228.2858
$ws.Range("M74").Value = -15.87176999999997
$ws.Range("N74").Value = -2976.2858
# Row 77
$ws.Range("H77").Value = 941.36957
$ws.Range("I77").Value = 889.87177
$ws.Range("J77").Value = 1228.2858
$ws.Range("K77").Value = 4449.35885
$ws.Range("L77").Value = 6141.429
$ws.Range("M77").Value = -81.35884999999962
$ws.Range("N77").Value = -14877.429
# Row 102
$ws.Range("H102").Value = 2330
$ws.Range("I102").Value = 2899.8
$ws.Range("J102").Value = 905.5
$ws.Range("K102").Value = 2899.8
$ws.Range("L102").Value = 905.5
$ws.Range("M102").Value = -1277.8
$ws.Range("N102").Value = -4149.5
# Row 132
$ws.Range("H132").Value = 3968.7778
$ws.Range("I132").Value = 3456
$ws.Range("J132").Value = 4115.2856
$ws.Range("K132").Value = 10368
$ws.Range("L132").Value = 12345.8568
$ws.Range("M132").Value = -7838
$ws.Range("N132").Value = -17405.8568

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2411.2
$ws.Range("I16").Value = 1066.3334
$ws.Range("J16").Value = 2987.5715
$ws.Range("K16").Value = 1066.3334
$ws.Range("L16").Value = 2987.5715
$ws.Range("M16").Value = -779.3334
$ws.Range("N16").Value = -3561.5715
# Row 31
$ws.Range("H31").Value = 2333.0193
$ws.Range("I31").Value = 1911.4286
$ws.Range("J31").Value = 4103.7
$ws.Range("K31").Value = 1911.4286
$ws.Range("L31").Value = 4103.7
$ws.Range("M31").Value = -1616.4286
$ws.Range("N31").Value = -4693.7
# Row 34
$ws.Range("H34").Value = 2333.0193
$ws.Range("I34").Value = 1911.4286
$ws.Range("J34").Value = 4103.7
$ws.Range("K34").Value = 1911.4286
$ws.Range("L34").Value = 4103.7
$ws.Range("M34").Value = -1709.4286
$ws.Range("N34").Value = -4507.7
# Row 113
$ws.Range("H113").Value = 2411.2
$ws.Range("I113").Value = 1066.3334
$ws.Range("J113").Value = 2987.5715
$ws.Range("K113").Value = 1066.3334
$ws.Range("L113").Value = 2987.5715
$ws.Range("M113").Value = 1103.6666
$ws.Range("N113").Value = -7327.5715
# Row 132
$ws.Range("H132").Value = 1589.2727
$ws.Range("I132").Value = 1075.3
$ws.Range("J132").Value = 2380
$ws.Range("K132").Value = 3225.9
$ws.Range("L132").Value = 7140
$ws.Range("M132").Value = -695.8999999999996
$ws.Range("N132").Value = -12200

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 84.833336
$ws.Range("I12").Value = 38.833332
$ws.Range("K12").Value = 116.499996
$ws.Range("M12").Value = 56.500004
# Row 97
$ws.Range("H97").Value = 433.0909
$ws.Range("I97").Value = 419.44446
$ws.Range("J97").Value = 494.5
$ws.Range("K97").Value = 1258.33338
$ws.Range("L97").Value = 1483.5
$ws.Range("M97").Value = -762.33338
$ws.Range("N97").Value = -2475.5

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 7131.3335
$ws.Range("I132").Value = 9333.333000000001
$ws.Range("J132").Value = 4929.3335
$ws.Range("K132").Value = 27999.999
$ws.Range("L132").Value = 14788.0005
$ws.Range("M132").Value = -25469.999
$ws.Range("N132").Value = -19848.0005

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 405938.97
$ws.Range("I40").Value = 533389.2
$ws.Range("K40").Value = 533389.2
$ws.Range("M40").Value = -533253.2
# Row 136
$ws.Range("H136").Value = 4814.6284
$ws.Range("I136").Value = 2674.4348
$ws.Range("J136").Value = 8916.666999999999
$ws.Range("K136").Value = 8023.3044
$ws.Range("L136").Value = 26750.001
$ws.Range("M136").Value = -5473.3044
$ws.Range("N136").Value = -31850.001

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 3321.652
$ws.Range("I96").Value = 2383.3333
$ws.Range("J96").Value = 3652.8235
$ws.Range("K96").Value = 2383.3333
$ws.Range("L96").Value = 3652.8235
$ws.Range("M96").Value = -1010.3333
$ws.Range("N96").Value = -6398.8235
# Row 122
$ws.Range("H122").Value = 1133.3438
$ws.Range("I122").Value = 1147.1724
$ws.Range("J122").Value = 999.6667
$ws.Range("K122").Value = 3441.5172
$ws.Range("L122").Value = 2999.0001
$ws.Range("M122").Value = -991.5171999999998
$ws.Range("N122").Value = -7899.0001
# Row 126
$ws.Range("H126").Value = 1033.9166
$ws.Range("I126").Value = 1020.7
$ws.Range("K126").Value = 3062.1
$ws.Range("M126").Value = -592.1000000000004
